$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the taxonomy/category values on the active "createNewCategory" sheet
$ws.Range("E2").Value = "AutomationTestTaxonomy"
$ws.Range("F2").Value = "AutomationTestCat1"

# Move the selection to F2 only
$ws.Range("F2").Select()
